# Insert one new data row at row 285 (this pushes the existing rows 285..322
# down to 286..323, which matches the target diff exactly) and then fill the
# brand-new row 285 with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 285..322 down by one row.
$ws.Rows.Item(285).Insert()

# Populate the newly inserted row 285 with the new record.
$ws.Cells.Item(285, 1).Value  = 4
$ws.Cells.Item(285, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(285, 3).Value  = "Los Lagos"
$ws.Cells.Item(285, 4).Value  = 44984
$ws.Cells.Item(285, 5).Value  = 10
$ws.Cells.Item(285, 6).Value  = 100112039
$ws.Cells.Item(285, 7).Value  = "Ciboulette"
$ws.Cells.Item(285, 8).Value  = "Sin especificar"
$ws.Cells.Item(285, 9).Value  = "Primera"
$ws.Cells.Item(285, 10).Value = 40
$ws.Cells.Item(285, 11).Value = 3500
$ws.Cells.Item(285, 12).Value = 3500
$ws.Cells.Item(285, 13).Value = 3500
$ws.Cells.Item(285, 14).Value = "$/docena de atados"
$ws.Cells.Item(285, 15).Value = "Región Metropolitana"
$ws.Cells.Item(285, 16).Value = 1167
$ws.Cells.Item(285, 17).Value = 3
$ws.Cells.Item(285, 18).Value = "Hortaliza"
